# Daily "cryptos" data refresh (GitHub Actions bot edit).
# Updates the Price (D) / Volume(1h) (E) columns for each coin row, plus a
# couple of rows whose rank (and therefore row position) changed, swapping
# NEARProtocol/Aptos (rows 25-26) and Stacks/USDe (rows 42-43).
#
# NOTE: several new Price values (e.g. "0.999", "218.52") are valid-looking
# numbers, but the sheet stores Price/Volume as literal text. A leading
# apostrophe forces Excel to keep the assigned string as text instead of
# silently coercing it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.022.49'
$ws.Range("E2").Value = '  +3.75%  '
$ws.Range("D3").Value = '3.221.22'
$ws.Range("E3").Value = '  +2.31%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '''218.52'
$ws.Range("E5").Value = '  +6.47%  '
$ws.Range("D6").Value = '''632.51'
$ws.Range("E6").Value = '  +4.74%  '
$ws.Range("D7").Value = '''0.395'
$ws.Range("E7").Value = '  +9.34%  '
$ws.Range("E8").Value = '  +7.45%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '3.218.36'
$ws.Range("E10").Value = '  +2.24%  '
$ws.Range("D11").Value = '''0.577'
$ws.Range("E11").Value = '  +8.87%  '
$ws.Range("D12").Value = '''0.180'
$ws.Range("E12").Value = '  +1.68%  '
$ws.Range("D13").Value = '''0.0000262'
$ws.Range("E13").Value = '  +9.18%  '
$ws.Range("D14").Value = '''5.44'
$ws.Range("E14").Value = '  +4.08%  '
$ws.Range("D15").Value = '''33.74'
$ws.Range("E15").Value = '  +6.33%  '
$ws.Range("D16").Value = '3.814.72'
$ws.Range("E16").Value = '  +2.09%  '
$ws.Range("D17").Value = '89.756.92'
$ws.Range("E17").Value = '  +3.68%  '
$ws.Range("D18").Value = '3.226.72'
$ws.Range("E18").Value = '  +1.69%  '
$ws.Range("D19").Value = '''0.0000244'
$ws.Range("E19").Value = '  +90.89%  '
$ws.Range("D20").Value = '''3.50'
$ws.Range("E20").Value = '  +19.38%  '
$ws.Range("D21").Value = '''13.63'
$ws.Range("E21").Value = '  +3.00%  '
$ws.Range("D22").Value = '''442.24'
$ws.Range("E22").Value = '  +8.26%  '
$ws.Range("D23").Value = '''8.70'
$ws.Range("E23").Value = '  +3.84%  '
$ws.Range("D24").Value = '''5.12'
$ws.Range("E24").Value = '  +2.59%  '
$ws.Range("B25").Value = 'NEARProtocol'
$ws.Range("C25").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D25").Value = '''5.28'
$ws.Range("E25").Value = '  +4.01%  '
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").Value = '''12.07'
$ws.Range("E26").Value = '  +5.26%  '
$ws.Range("D27").Value = '''82.81'
$ws.Range("E27").Value = '  +13.63%  '
$ws.Range("D28").Value = '3.395.76'
$ws.Range("E28").Value = '  +1.52%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("E30").Value = '  +0.94%  '
$ws.Range("D31").Value = '''0.999'
$ws.Range("E31").Value = '  +0.65%  '
$ws.Range("D32").Value = '''4.22'
$ws.Range("E32").Value = '  +42.28%  '
$ws.Range("D33").Value = '''8.62'
$ws.Range("E33").Value = '  +5.21%  '
$ws.Range("D34").Value = '''548.19'
$ws.Range("E34").Value = '  +2.73%  '
$ws.Range("D35").Value = '''7.08'
$ws.Range("E35").Value = '  +8.82%  '
$ws.Range("E36").Value = '  +4.82%  '
$ws.Range("D37").Value = '''1.32'
$ws.Range("E37").Value = '  +4.73%  '
$ws.Range("D38").Value = '''22.48'
$ws.Range("E38").Value = '  +4.56%  '
$ws.Range("D39").Value = '''22.40'
$ws.Range("E39").Value = '  +2.85%  '
$ws.Range("E40").Value = '  -0.96%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''1.95'
$ws.Range("E42").Value = '  +4.08%  '
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '''0.378'
$ws.Range("E44").Value = '  +2.14%  '
$ws.Range("D45").Value = '''147.09'
$ws.Range("E45").Value = '  -1.82%  '
$ws.Range("D46").Value = '''174.36'
$ws.Range("E46").Value = '  +2.33%  '
$ws.Range("D47").Value = '''43.93'
$ws.Range("E47").Value = '  +2.22%  '
$ws.Range("D48").Value = '''0.759'
$ws.Range("E48").Value = '  +11.20%  '
$ws.Range("D49").Value = '''1.27'
$ws.Range("E49").Value = '  +2.52%  '
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("D51").Value = '''0.628'
$ws.Range("E51").Value = '  +8.08%  '
